# Workout guru installation.docx — add the Eclipse.ini / Weka notes
# right after the "Goto Android tab ..." paragraph and before the
# "Algorithm considerations:" paragraph.

$d = $word.ActiveDocument

# Locate the paragraph that ends the Eclipse/Add-library instructions.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Goto Android tab*") {
        $anchor = $p
    }
}

# Insert point: right before the anchor paragraph's own end-of-paragraph
# mark, i.e. at the very end of the anchor paragraph's text.
$insertPos = $anchor.Range.End - 1
$r = $d.Range($insertPos, $insertPos)

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml = @"
<w:p $wns>
  <w:r><w:t>Eclipse.</w:t></w:r>
  <w:r><w:t xml:space="preserve">ini file change: </w:t></w:r>
  <w:r><w:t>on mac right click on eclipse.app and click package contents. That will lead to the in file</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Times" w:eastAsia="Times New Roman" w:hAnsi="Times" w:cs="Times New Roman"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Times New Roman"/>
      <w:color w:val="444444"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FAFAFA"/>
    </w:rPr>
    <w:t>--launcher.XXMaxPermSize1024m -vmargs -Xms256m -Xmx1024m</w:t>
  </w:r>
</w:p>
<w:p $wns/>
<w:p $wns>
  <w:r><w:t>Use Weka for android jar instead of regular weka.jar</w:t></w:r>
</w:p>
<w:p $wns/>
"@

$r.InsertXML($xml) | Out-Null
